{"js": "// Update the two-digit-divided-by-one-digit practice table: each populated\n// cell's \"NN\u00f7N=\" expression is replaced by a new one, in document order.\n// (Several old expressions repeat, e.g. \"57\u00f78=\" appears twice with two\n// different replacements, so we walk the table in order rather than doing a\n// single blind find/replace.)\n\nconst newValues = [\n  \"94\u00f73=\",\n  \"33\u00f74=\",\n  \"30\u00f79=\",\n  \"78\u00f78=\",\n  \"94\u00f73=\",\n  \"21\u00f75=\",\n  \"74\u00f74=\",\n  \"61\u00f74=\",\n  \"19\u00f79=\",\n  \"62\u00f77=\",\n  \"50\u00f76=\",\n  \"76\u00f78=\",\n  \"73\u00f74=\",\n  \"89\u00f72=\",\n  \"63\u00f77=\",\n  \"54\u00f77=\",\n  \"57\u00f79=\",\n  \"93\u00f77=\",\n  \"80\u00f78=\",\n  \"17\u00f79=\",\n  \"94\u00f77=\",\n  \"13\u00f75=\",\n  \"68\u00f79=\",\n  \"74\u00f73=\",\n  \"69\u00f72=\",\n];\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\ntable.load(\"values\");\nawait context.sync();\n\nconst rows = table.values.length;\nconst cols = table.values[0].length;\n\nlet i = 0;\nfor (let r = 0; r < rows; r++) {\n  for (let c = 0; c < cols; c++) {\n    const text = table.values[r][c];\n    if (text !== \"\") {\n      table.getCell(r, c).value = newValues[i];\n      i++;\n    }\n  }\n}\n\nawait context.sync();\n", "ps1": "# Update the two-digit-divided-by-one-digit practice table: each populated\n# cell's \"NN\u00f7N=\" expression is replaced by a new one, in document order.\n# (Several old expressions repeat, e.g. \"57\u00f78=\" appears twice with two\n# different replacements, so we walk the table in order rather than doing a\n# single blind text Find/Replace.)\n\n$newValues = @(\n    \"94\u00f73=\",\n    \"33\u00f74=\",\n    \"30\u00f79=\",\n    \"78\u00f78=\",\n    \"94\u00f73=\",\n    \"21\u00f75=\",\n    \"74\u00f74=\",\n    \"61\u00f74=\",\n    \"19\u00f79=\",\n    \"62\u00f77=\",\n    \"50\u00f76=\",\n    \"76\u00f78=\",\n    \"73\u00f74=\",\n    \"89\u00f72=\",\n    \"63\u00f77=\",\n    \"54\u00f77=\",\n    \"57\u00f79=\",\n    \"93\u00f77=\",\n    \"80\u00f78=\",\n    \"17\u00f79=\",\n    \"94\u00f77=\",\n    \"13\u00f75=\",\n    \"68\u00f79=\",\n    \"74\u00f73=\",\n    \"69\u00f72=\"\n)\n\n$d = $word.ActiveDocument\n$tbl = $d.Tables.Item(1)\n\n$i = 0\nfor ($r = 1; $r -le $tbl.Rows.Count; $r++) {\n    for ($c = 1; $c -le $tbl.Columns.Count; $c++) {\n        $cell = $tbl.Cell($r, $c)\n        $cellText = $cell.Range.Text\n        # Strip the trailing end-of-cell marker(s) before checking content.\n        $plain = $cellText -replace \"[\\x07\\r]+$\", \"\"\n        if ($plain -ne \"\") {\n            $cell.Range.Text = $newValues[$i]\n            $i = $i + 1\n        }\n    }\n}\n"}
